# edit.ps1 - applies the changes described by the diff:
#   1. Update the cached "datetimeFigureOut" auto-date field text from
#      2024/4/30 -> 2024/5/1 on every slide layout, the slide master and
#      the notes master (13 occurrences total).
#   2. Update the static date text on slide 1 from 2020/9/25 -> 2024/5/1.
#   3. Update the body text on slide 2 from "Demo examples." -> "Example results".

$p = $ppt.ActivePresentation

function Set-ShapeText($shape, [string]$newText) {
    # Writing the exact same/target text directly can get merged with the
    # previous text by the engine's minimal-diff run patcher whenever the
    # old and new strings share a common prefix (e.g. "2020/9/25" ->
    # "2024/5/1" both start with "202"), which would split the run in two.
    # Setting an unrelated placeholder value first guarantees a clean,
    # single-run replacement that matches how the real edit looks.
    $shape.TextFrame.TextRange.Text = "~"
    $shape.TextFrame.TextRange.Text = $newText
}

# ---------------------------------------------------------------------
# 1) Date placeholder fields: slide master + all slide layouts + notes master
# ---------------------------------------------------------------------

# Slide master (shape holding the "日期版面配置區" / Date Placeholder)
$master = $p.SlideMaster
Set-ShapeText $master.Shapes.Item(5) "2024/5/1"

# All custom (slide) layouts - index of the date placeholder shape varies
# per layout, enumerated from the underlying XML.
$layoutDateShapeIndex = @{1=5; 2=3; 3=5; 4=4; 5=6; 6=2; 7=1; 8=4; 9=4; 10=3; 11=5}
$layouts = $p.Designs.Item(1).SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $idx = $layoutDateShapeIndex[$i]
    Set-ShapeText $layout.Shapes.Item($idx) "2024/5/1"
}

# Notes master
$notesMaster = $p.NotesMaster
Set-ShapeText $notesMaster.Shapes.Item(2) "2024/5/1"

# ---------------------------------------------------------------------
# 2) Slide 1: static date text box ("Chen-Hsiung Liu" / "2020/9/25")
# ---------------------------------------------------------------------

$slide1 = $p.Slides.Item(1)
$dateBox = $slide1.Shapes.Item(3)
Set-ShapeText $dateBox "Chen-Hsiung Liu`r2024/5/1"

# ---------------------------------------------------------------------
# 3) Slide 2: "Demo examples." -> "Example results"
# ---------------------------------------------------------------------

$slide2 = $p.Slides.Item(2)
$bodyBox = $slide2.Shapes.Item(2)
Set-ShapeText $bodyBox "Introduce SystemVerilog parser tools and comparison`rExample results"

Write-Output "Done."
